$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 730.4
$ws.Range("I8").Value = 200.8
$ws.Range("J8").Value = 1260
$ws.Range("K8").Value = 602.4000000000001
$ws.Range("L8").Value = 3780
$ws.Range("M8").Value = -463.4000000000001
$ws.Range("N8").Value = -4058

$ws.Range("H116").Value = 4943.1
$ws.Range("I116").Value = 3903.9375
$ws.Range("J116").Value = 9099.75
$ws.Range("K116").Value = 3903.9375
$ws.Range("L116").Value = 9099.75
$ws.Range("M116").Value = -461.9375
$ws.Range("N116").Value = -15983.75

$ws.Range("H138").Value = 5212738.5
$ws.Range("I138").Value = 2083.2222
$ws.Range("J138").Value = 6415197.5
$ws.Range("K138").Value = 6249.6666
$ws.Range("L138").Value = 19245592.5
$ws.Range("M138").Value = -1109.6666
$ws.Range("N138").Value = -19255872.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2666.6667
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H45").Value = 18739.334
$ws.Range("J45").Value = 2492.25
$ws.Range("L45").Value = 2492.25
$ws.Range("N45").Value = -3246.25

$ws.Range("H74").Value = 5054.8945
$ws.Range("I74").Value = 2237.2354
$ws.Range("J74").Value = 29005
$ws.Range("K74").Value = 2237.2354
$ws.Range("L74").Value = 29005
$ws.Range("M74").Value = -1363.2354
$ws.Range("N74").Value = -30753

$ws.Range("H76").Value = 15999
$ws.Range("J76").Value = 15999
$ws.Range("L76").Value = 15999
$ws.Range("N76").Value = -16675

$ws.Range("H77").Value = 5054.8945
$ws.Range("I77").Value = 2237.2354
$ws.Range("J77").Value = 29005
$ws.Range("K77").Value = 11186.177
$ws.Range("L77").Value = 145025
$ws.Range("M77").Value = -6818.177
$ws.Range("N77").Value = -153761

$ws.Range("H79").Value = 15999
$ws.Range("J79").Value = 15999
$ws.Range("L79").Value = 15999
$ws.Range("N79").Value = -18339

$ws.Range("H86").Value = 65999
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 65999
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H110").Value = 7163.174
$ws.Range("I110").Value = 10196.818
$ws.Range("J110").Value = 4382.3335
$ws.Range("K110").Value = 10196.818
$ws.Range("L110").Value = 4382.3335
$ws.Range("M110").Value = -8151.817999999999
$ws.Range("N110").Value = -8472.333500000001

$ws.Range("H116").Value = 2666.6667
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H122").Value = 3077.1538
$ws.Range("I122").Value = 2395.5557
$ws.Range("K122").Value = 7186.6671
$ws.Range("M122").Value = -4736.6671

$ws.Range("H132").Value = 3046.205
$ws.Range("I132").Value = 2717.1936
$ws.Range("K132").Value = 8151.5808
$ws.Range("M132").Value = -5621.5808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2666.6667
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H12").Value = 447.5
$ws.Range("I12").Value = 295
$ws.Range("K12").Value = 295
$ws.Range("M12").Value = -127

$ws.Range("H20").Value = 4284.636
$ws.Range("I20").Value = 5146.4
$ws.Range("J20").Value = 3566.5
$ws.Range("K20").Value = 5146.4
$ws.Range("L20").Value = 3566.5
$ws.Range("M20").Value = -4899.4
$ws.Range("N20").Value = -4060.5

$ws.Range("H86").Value = 8359
$ws.Range("I86").Value = 10602
$ws.Range("J86").Value = 4994.5
$ws.Range("K86").Value = 10602
$ws.Range("L86").Value = 4994.5
$ws.Range("M86").Value = -9479
$ws.Range("N86").Value = -7240.5

$ws.Range("H89").Value = 8359
$ws.Range("I89").Value = 10602
$ws.Range("J89").Value = 4994.5
$ws.Range("K89").Value = 53010
$ws.Range("L89").Value = 24972.5
$ws.Range("M89").Value = -47394
$ws.Range("N89").Value = -36204.5

$ws.Range("H99").Value = 5111.5713
$ws.Range("J99").Value = 7088.9
$ws.Range("L99").Value = 7088.9
$ws.Range("N99").Value = -10084.9

$ws.Range("H105").Value = 6815.778
$ws.Range("I105").Value = 6917.75
$ws.Range("K105").Value = 6917.75
$ws.Range("M105").Value = -5170.75

$ws.Range("H107").Value = 2766.923
$ws.Range("I107").Value = 2766.923
$ws.Range("K107").Value = 2766.923
$ws.Range("M107").Value = -846.9229999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3653.5557
$ws.Range("I31").Value = 3251.5
$ws.Range("J31").Value = 3975.2
$ws.Range("K31").Value = 3251.5
$ws.Range("L31").Value = 3975.2
$ws.Range("M31").Value = -2956.5
$ws.Range("N31").Value = -4565.2

$ws.Range("H34").Value = 3653.5557
$ws.Range("I34").Value = 3251.5
$ws.Range("J34").Value = 3975.2
$ws.Range("K34").Value = 3251.5
$ws.Range("L34").Value = 3975.2
$ws.Range("M34").Value = -3049.5
$ws.Range("N34").Value = -4379.2

$ws.Range("H39").Value = 17500
$ws.Range("I39").Value = 15000
$ws.Range("J39").Value = 20000
$ws.Range("K39").Value = 15000
$ws.Range("L39").Value = 20000
$ws.Range("M39").Value = -14609
$ws.Range("N39").Value = -20782

$ws.Range("H49").Value = 17500
$ws.Range("I49").Value = 15000
$ws.Range("J49").Value = 20000
$ws.Range("K49").Value = 15000
$ws.Range("L49").Value = 20000
$ws.Range("M49").Value = -14818
$ws.Range("N49").Value = -20364

$ws.Range("H105").Value = 2092.3572
$ws.Range("J105").Value = 2177.75
$ws.Range("L105").Value = 2177.75
$ws.Range("N105").Value = -5671.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 8441.429
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 8441.429
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 25324.287
$ws.Range("N80").Value = -27196.287
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 8441.429
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 8441.429
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 75972.861
$ws.Range("N83").Value = -85332.861
$ws.Range("M83").ClearContents()

$ws.Range("H113").Value = 1326
$ws.Range("J113").Value = 1290
$ws.Range("L113").Value = 3870
$ws.Range("N113").Value = -8210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 25748.375
$ws.Range("I70").Value = 25747.5
$ws.Range("K70").Value = 25747.5
$ws.Range("M70").Value = -25477.5

$ws.Range("H73").Value = 25748.375
$ws.Range("I73").Value = 25747.5
$ws.Range("K73").Value = 25747.5
$ws.Range("M73").Value = -24811.5

$ws.Range("H97").Value = 2177.75
$ws.Range("I97").Value = 2403.8333
$ws.Range("J97").Value = 1499.5
$ws.Range("K97").Value = 2403.8333
$ws.Range("L97").Value = 1499.5
$ws.Range("M97").Value = -1907.8333
$ws.Range("N97").Value = -2491.5

$ws.Range("H113").Value = 1621.3334
$ws.Range("I113").Value = 1663.1428
$ws.Range("K113").Value = 1663.1428
$ws.Range("M113").Value = 506.8571999999999

$ws.Range("H122").Value = 1596.5333
$ws.Range("I122").Value = 1596.5333
$ws.Range("K122").Value = 4789.5999
$ws.Range("M122").Value = -2339.5999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 25000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H61").Value = 2002.8889
$ws.Range("I61").Value = 1426.4103
$ws.Range("J61").Value = 5750
$ws.Range("K61").Value = 1426.4103
$ws.Range("L61").Value = 5750
$ws.Range("M61").Value = -1224.4103
$ws.Range("N61").Value = -6154

$ws.Range("H82").Value = 7392.75
$ws.Range("I82").Value = 10446.818
$ws.Range("J82").Value = 3660
$ws.Range("K82").Value = 10446.818
$ws.Range("L82").Value = 3660
$ws.Range("M82").Value = -10085.818
$ws.Range("N82").Value = -4382

$ws.Range("H85").Value = 7392.75
$ws.Range("I85").Value = 10446.818
$ws.Range("J85").Value = 3660
$ws.Range("K85").Value = 10446.818
$ws.Range("L85").Value = 3660
$ws.Range("M85").Value = -9198.817999999999
$ws.Range("N85").Value = -6156

$ws.Range("H113").Value = 2002.8889
$ws.Range("I113").Value = 1426.4103
$ws.Range("J113").Value = 5750
$ws.Range("K113").Value = 1426.4103
$ws.Range("L113").Value = 5750
$ws.Range("M113").Value = 743.5897
$ws.Range("N113").Value = -10090

$ws.Range("H132").Value = 2699.353
$ws.Range("I132").Value = 2023.2222
$ws.Range("K132").Value = 6069.6666
$ws.Range("M132").Value = -3539.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 4004505.2
$ws.Range("I22").Value = 6666675.5
$ws.Range("J22").Value = 11250
$ws.Range("K22").Value = 6666675.5
$ws.Range("L22").Value = 11250
$ws.Range("M22").Value = -6666382.5
$ws.Range("N22").Value = -11836

$ws.Range("H107").Value = 2242.7144
$ws.Range("I107").Value = 2666.3333
$ws.Range("K107").Value = 7998.999899999999
$ws.Range("M107").Value = -6078.999899999999

$ws.Range("H113").Value = 1333.3462
$ws.Range("I113").Value = 1418.55
$ws.Range("J113").Value = 1049.3334
$ws.Range("K113").Value = 4255.65
$ws.Range("L113").Value = 3148.0002
$ws.Range("M113").Value = -2085.65
$ws.Range("N113").Value = -7488.0002

$ws.Range("H126").Value = 3378.5
$ws.Range("J126").Value = 3575.2856
$ws.Range("L126").Value = 10725.8568
$ws.Range("N126").Value = -15665.8568

$ws.Range("H136").Value = 2045.9
$ws.Range("I136").Value = 1807.1471
$ws.Range("K136").Value = 5421.4413
$ws.Range("M136").Value = -2871.4413
